$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# --- Header row (row 1): add date / legislator_name / legislator_id columns ---
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "date"

$ws.Range("G1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "legislator_name"

$ws.Range("G1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "legislator_id"

# --- Data rows (2-13): fill the new columns for every stock entry ---
# date column (H) - must stay plain text "2012-03-30", not be auto-converted to a date
$ws.Range("H2:H13").NumberFormat = "@"
$ws.Range("H2:H13").Value = "2012-03-30"
$ws.Range("F2:F13").Copy()
$ws.Range("H2:H13").PasteSpecial(-4122)

# legislator_name column (I)
$ws.Range("G2:G13").Copy()
$ws.Range("I2:I13").PasteSpecial(-4122)
$ws.Range("I2:I13").Value = "黃志雄"

# legislator_id column (J) - numeric id
$ws.Range("G2:G13").Copy()
$ws.Range("J2:J13").PasteSpecial(-4122)
$ws.Range("J2:J13").Value = 1366

$excel.CutCopyMode = 0
